$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 so existing rows (Kraken2_0.0 ... Mothur)
# shift down by one. The old "NBC" row (formerly row 10, now row 11) is
# renamed to "CustomNBC" and its whole row is moved up into the new row 4,
# then the now-duplicate row 11 is deleted.

$ws.Rows.Item(4).Insert()

$srcRow = 11
$destRow = 4

$ws.Cells.Item($destRow, 1).Value = "CustomNBC"
for ($col = 2; $col -le 6; $col++) {
    $ws.Cells.Item($destRow, $col).Value = $ws.Cells.Item($srcRow, $col).Value2
}

$ws.Rows.Item($srcRow).Delete()
